$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from M1:P1 (style s=1, bordered/bold/centered header) to Q1:T1
$ws.Range("M1:P1").Copy()
$ws.Range("Q1:T1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row labels
$ws.Range("Q1").Value = "Estación más cercana 6"
$ws.Range("R1").Value = "Estación más cercana 7"
$ws.Range("S1").Value = "Inicio estación más cercana 6"
$ws.Range("T1").Value = "Inicio estación más cercana 7"

# Data rows
$ws.Cells.Item(2, 17).Value = "AC05"
$ws.Cells.Item(2, 18).Value = "CO10"
$ws.Cells.Item(2, 19).Value = "2014-03-04T02:50:02"
$ws.Cells.Item(2, 20).Value = "2014-03-04T02:50:03"
$ws.Cells.Item(3, 17).Value = "PX02"
$ws.Cells.Item(3, 18).Value = "PB11"
$ws.Cells.Item(3, 19).Value = "2014-04-04T07:57:49"
$ws.Cells.Item(3, 20).Value = "2014-04-04T07:57:50"
$ws.Cells.Item(4, 17).Value = "PB11"
$ws.Cells.Item(4, 18).Value = "PX03"
$ws.Cells.Item(4, 19).Value = "2014-04-12T22:35:21"
$ws.Cells.Item(4, 20).Value = "2014-04-12T22:35:22"
$ws.Cells.Item(5, 17).Value = "PB12"
$ws.Cells.Item(5, 18).Value = "PB11"
$ws.Cells.Item(5, 19).Value = "2014-03-24T12:22:36"
$ws.Cells.Item(5, 20).Value = "2014-03-24T12:22:37"
$ws.Cells.Item(6, 17).Value = "PB11"
$ws.Cells.Item(6, 18).Value = "PB12"
$ws.Cells.Item(6, 19).Value = "2014-03-24T05:22:54"
$ws.Cells.Item(6, 20).Value = "2014-03-24T05:22:56"
$ws.Cells.Item(7, 17).Value = "TA01"
$ws.Cells.Item(7, 18).Value = "HMBCX"
$ws.Cells.Item(7, 19).Value = "2014-03-15T09:31:36"
$ws.Cells.Item(7, 20).Value = "2014-03-15T09:31:36"
$ws.Cells.Item(8, 17).Value = "ROC1"
$ws.Cells.Item(8, 18).Value = "MT07"
$ws.Cells.Item(8, 19).Value = "2014-05-15T16:21:07"
$ws.Cells.Item(8, 20).Value = "2014-05-15T16:21:07"
$ws.Cells.Item(9, 17).Value = "PB11"
$ws.Cells.Item(9, 18).Value = "PX02"
$ws.Cells.Item(9, 19).Value = "2014-04-20T04:23:24"
$ws.Cells.Item(9, 20).Value = "2014-04-20T04:23:24"
$ws.Cells.Item(10, 17).Value = "PX03"
$ws.Cells.Item(10, 18).Value = "PB11"
$ws.Cells.Item(10, 19).Value = "2014-04-21T07:10:54"
$ws.Cells.Item(10, 20).Value = "2014-04-21T07:10:55"
$ws.Cells.Item(11, 17).Value = "PX03"
$ws.Cells.Item(11, 18).Value = "PB11"
$ws.Cells.Item(11, 19).Value = "2014-06-19T00:31:36"
$ws.Cells.Item(11, 20).Value = "2014-06-19T00:31:37"
$ws.Cells.Item(12, 17).Value = "MT07"
$ws.Cells.Item(12, 18).Value = "MT02"
$ws.Cells.Item(12, 19).Value = "2014-06-24T11:23:42"
$ws.Cells.Item(12, 20).Value = "2014-06-24T11:23:44"
$ws.Cells.Item(13, 17).Value = "MT07"
$ws.Cells.Item(13, 18).Value = "V25A"
$ws.Cells.Item(13, 19).Value = "2014-12-15T00:32:33"
$ws.Cells.Item(13, 20).Value = "2014-12-15T00:32:34"
$ws.Cells.Item(14, 17).Value = "PATCX"
$ws.Cells.Item(14, 18).Value = "PB11"
$ws.Cells.Item(14, 19).Value = "2015-02-14T14:47:48"
$ws.Cells.Item(14, 20).Value = "2015-02-14T14:47:48"
$ws.Cells.Item(15, 17).Value = "AC01"
$ws.Cells.Item(15, 18).Value = "AC05"
$ws.Cells.Item(15, 19).Value = "2015-03-02T14:02:24"
$ws.Cells.Item(15, 20).Value = "2015-03-02T14:02:33"
$ws.Cells.Item(16, 17).Value = "VA01"
$ws.Cells.Item(16, 18).Value = "ROC1"
$ws.Cells.Item(16, 19).Value = "2015-09-21T19:08:22"
$ws.Cells.Item(16, 20).Value = "2015-09-21T19:08:24"
$ws.Cells.Item(17, 17).Value = "GO04"
$ws.Cells.Item(17, 18).Value = "TLL"
$ws.Cells.Item(17, 19).Value = "2015-10-03T18:28:04"
$ws.Cells.Item(17, 20).Value = "2015-10-03T18:28:04"
$ws.Cells.Item(18, 17).Value = "CO02"
$ws.Cells.Item(18, 18).Value = "CO03"
$ws.Cells.Item(18, 19).Value = "2015-10-09T22:18:28"
$ws.Cells.Item(18, 20).Value = "2015-10-09T22:18:28"
$ws.Cells.Item(19, 17).Value = "CO05"
$ws.Cells.Item(19, 18).Value = "GO04"
$ws.Cells.Item(19, 19).Value = "2015-09-20T23:00:24"
$ws.Cells.Item(19, 20).Value = "2015-09-20T23:00:25"
$ws.Cells.Item(20, 17).Value = "MT07"
$ws.Cells.Item(20, 18).Value = "MT02"
$ws.Cells.Item(20, 19).Value = "2015-09-19T21:27:43"
$ws.Cells.Item(20, 20).Value = "2015-09-19T21:27:44"
$ws.Cells.Item(21, 17).Value = "IN40"
$ws.Cells.Item(21, 18).Value = "IN41"
$ws.Cells.Item(21, 19).Value = "2015-09-19T22:35:41"
$ws.Cells.Item(21, 20).Value = "2015-09-19T22:35:41"
$ws.Cells.Item(22, 17).Value = "TLL"
$ws.Cells.Item(22, 18).Value = "CO10"
$ws.Cells.Item(22, 19).Value = "2015-09-20T14:45:58"
$ws.Cells.Item(22, 20).Value = "2015-09-20T14:45:58"
$ws.Cells.Item(23, 17).Value = "TLL"
$ws.Cells.Item(23, 18).Value = "CO02"
$ws.Cells.Item(23, 19).Value = "2015-10-08T02:54:32"
$ws.Cells.Item(23, 20).Value = "2015-10-08T02:54:34"
$ws.Cells.Item(24, 17).Value = "CO03"
$ws.Cells.Item(24, 18).Value = "ROC1"
$ws.Cells.Item(24, 19).Value = "2015-12-13T01:17:59"
$ws.Cells.Item(24, 20).Value = "2015-12-13T01:18:00"
$ws.Cells.Item(25, 17).Value = "TLL"
$ws.Cells.Item(25, 18).Value = "CO02"
$ws.Cells.Item(25, 19).Value = "2016-02-16T15:49:49"
$ws.Cells.Item(25, 20).Value = "2016-02-16T15:49:50"
$ws.Cells.Item(26, 17).Value = "TLL"
$ws.Cells.Item(26, 18).Value = "CO02"
$ws.Cells.Item(26, 19).Value = "2016-02-22T18:47:12"
$ws.Cells.Item(26, 20).Value = "2016-02-22T18:47:13"
$ws.Cells.Item(27, 17).Value = "PB01"
$ws.Cells.Item(27, 18).Value = "PB07"
$ws.Cells.Item(27, 19).Value = "2016-03-04T07:58:31"
$ws.Cells.Item(27, 20).Value = "2016-03-04T07:58:31"
$ws.Cells.Item(28, 17).Value = "IN40"
$ws.Cells.Item(28, 18).Value = "GO04"
$ws.Cells.Item(28, 19).Value = "2016-02-25T11:49:25"
$ws.Cells.Item(28, 20).Value = "2016-02-25T11:49:26"
$ws.Cells.Item(29, 17).Value = "PB11"
$ws.Cells.Item(29, 18).Value = "PB13"
$ws.Cells.Item(29, 19).Value = "2016-06-03T13:47:05"
$ws.Cells.Item(29, 20).Value = "2016-06-03T13:47:07"
$ws.Cells.Item(30, 17).Value = "TLL"
$ws.Cells.Item(30, 18).Value = "AC04"
$ws.Cells.Item(30, 19).Value = "2016-06-05T09:31:10"
$ws.Cells.Item(30, 20).Value = "2016-06-05T09:31:10"
$ws.Cells.Item(31, 17).Value = "CO03"
$ws.Cells.Item(31, 18).Value = "V25A"
$ws.Cells.Item(31, 19).Value = "2016-06-25T10:47:52"
$ws.Cells.Item(31, 20).Value = "2016-06-25T10:47:52"
$ws.Cells.Item(32, 17).Value = "TLL"
$ws.Cells.Item(32, 18).Value = "CO02"
$ws.Cells.Item(32, 19).Value = "2016-07-05T09:53:19"
$ws.Cells.Item(32, 20).Value = "2016-07-05T09:53:21"
$ws.Cells.Item(33, 17).Value = "PB06"
$ws.Cells.Item(33, 18).Value = "PB03"
$ws.Cells.Item(33, 19).Value = "2014-03-18T20:17:44"
$ws.Cells.Item(33, 20).Value = "2014-03-18T20:17:49"
$ws.Cells.Item(34, 17).Value = "MT07"
$ws.Cells.Item(34, 18).Value = "MT02"
$ws.Cells.Item(34, 19).Value = "2017-10-07T10:37:36"
$ws.Cells.Item(34, 20).Value = "2017-10-07T10:37:37"
$ws.Cells.Item(35, 17).Value = "PB12"
$ws.Cells.Item(35, 18).Value = "PATCX"
$ws.Cells.Item(35, 19).Value = "2015-07-23T00:27:56"
$ws.Cells.Item(35, 20).Value = "2015-07-23T00:27:57"
$ws.Cells.Item(36, 17).Value = "IN40"
$ws.Cells.Item(36, 18).Value = "IN41"
$ws.Cells.Item(36, 19).Value = "2018-06-29T22:17:44"
$ws.Cells.Item(36, 20).Value = "2018-06-29T22:17:44"
$ws.Cells.Item(37, 17).Value = "PB19"
$ws.Cells.Item(37, 18).Value = "PB05"
$ws.Cells.Item(37, 19).Value = "2019-03-28T21:01:14"
$ws.Cells.Item(37, 20).Value = "2019-03-28T21:01:18"
$ws.Cells.Item(38, 17).Value = "MT01"
$ws.Cells.Item(38, 18).Value = "BO01"
$ws.Cells.Item(38, 19).Value = "2019-08-02T00:09:27"
$ws.Cells.Item(38, 20).Value = "2019-08-02T00:09:29"
$ws.Cells.Item(39, 17).Value = "PATCX"
$ws.Cells.Item(39, 18).Value = "PB11"
$ws.Cells.Item(39, 19).Value = "2020-01-09T16:26:15"
$ws.Cells.Item(39, 20).Value = "2020-01-09T16:26:15"
$ws.Cells.Item(40, 17).Value = "MT07"
$ws.Cells.Item(40, 18).Value = "MT02"
$ws.Cells.Item(40, 19).Value = "2020-10-06T05:16:38"
$ws.Cells.Item(40, 20).Value = "2020-10-06T05:16:40"
$ws.Cells.Item(41, 17).Value = "GO04"
$ws.Cells.Item(41, 18).Value = "TLL"
$ws.Cells.Item(41, 19).Value = "2020-11-12T10:24:56"
$ws.Cells.Item(41, 20).Value = "2020-11-12T10:24:56"
